$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF with header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the existing header formatting (bold, centered, thin border) onto the
# new header cells so they reuse the same style as the other header cells.
$ws.Range("AA1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-56: Wins=69, Losses=93, Ties=0 (team record repeated per player row)
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 69   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 93   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
